# Applies the "output generated at 456a3b4" update to 杭州-漫展信息.xlsx
#
# Summary of the change:
#  - 展览   (sheet 1, "Exhibitions"): "want to go" counts (column F) bumped for many rows.
#  - 演出   (sheet 2, "Performances"): the cancelled 2024-07-25 show (row 2) is removed,
#           every later row shifts up by one, the running index in column A is
#           renumbered, and a handful of rows get updated "want to go" counts (column F).
#  - 本地生活 (sheet 3, "Local life"): "want to go" counts (column F) bumped.
#  - 全部类型 (sheet 4, "All types"): "want to go" counts (column F) bumped.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: 展览 - update column F ("想去人数") for the listed rows.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    2  = 864
    3  = 1795
    4  = 81
    5  = 545
    7  = 1380
    8  = 2113
    11 = 2428
    12 = 667
    14 = 4025
    16 = 375
    17 = 3139
    18 = 855
    21 = 175
    22 = 2088
    23 = 1195
    24 = 3
    25 = 1973
    26 = 399
    27 = 216
    28 = 26
    29 = 8740
    30 = 5793
    31 = 358
    32 = 185
    33 = 770
    34 = 11
    35 = 778
    36 = 3478
    40 = 51
    41 = 198
    42 = 170
    43 = 4665
    45 = 887
    46 = 85
    47 = 414
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

# ---------------------------------------------------------------------------
# Sheet 2: 演出 - drop the cancelled "07-25" show in row 2 (everything below
# shifts up by one row), then renumber column A (the running index) and fix
# up the "want to go" counts (column F) that differ from a plain shift.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Rows.Item(2).Delete() | Out-Null

# Renumber the running index in column A (0 for the header, 1..25 for data).
$lastRow = $ws2.Cells.Item($ws2.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws2.Cells.Item($r, 1).Value = $r - 1
}

$sheet2Updates = @{
    2  = 88
    14 = 110
    15 = 27
    17 = 3414
}
foreach ($row in $sheet2Updates.Keys) {
    $ws2.Cells.Item($row, 6).Value = $sheet2Updates[$row]
}

# ---------------------------------------------------------------------------
# Sheet 3: 本地生活 - update column F ("想去人数") for the listed rows.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$sheet3Updates = @{
    2 = 8402
    3 = 374
    4 = 1338
}
foreach ($row in $sheet3Updates.Keys) {
    $ws3.Cells.Item($row, 6).Value = $sheet3Updates[$row]
}

# ---------------------------------------------------------------------------
# Sheet 4: 全部类型 - update column F ("想去人数") for the listed rows.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    2  = 8402
    3  = 864
    4  = 374
    5  = 1338
    6  = 88
    7  = 1795
    8  = 81
    9  = 545
    10 = 1380
    11 = 2113
    15 = 4025
    16 = 375
    17 = 3139
    18 = 855
    21 = 175
    22 = 2088
    26 = 1195
    28 = 1973
    29 = 110
    30 = 399
    31 = 216
    32 = 26
    33 = 8740
    35 = 358
    36 = 778
    39 = 51
    40 = 198
    42 = 170
    43 = 888
    44 = 85
    45 = 414
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
